# Upload new version with timestamp
# Fills in the first data row of the report with a sold item and refreshes
# the "generated at" timestamp shown at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row counter for the first data row.
$ws.Range("A7").Value = 1

# Item name (C7:G7 is merged) - mark as Text so it is stored verbatim.
$ws.Range("C7:G7").NumberFormat = "@"
$ws.Range("C7").Value = "FAYCID HAIR OIL 60 ML"

# Time column (H7:K7 is merged) - "8:0" must stay literal text, not a time value.
$ws.Range("H7:K7").NumberFormat = "@"
$ws.Range("H7").Value = "8:0"

# Quantity (L7:M7 is merged) - keep the cell's original number format (#,##0.##...)
# but store the value as literal text "1", matching the source report export.
$ws.Range("L7:M7").NumberFormat = "@"
$ws.Range("L7").Value = "1"
$ws.Range("L7:M7").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"

# Unit price (N7:O7 is merged, shares formatting with the name column) - literal text.
$ws.Range("N7:O7").NumberFormat = "@"
$ws.Range("N7").Value = "30.00"

# Total price (P7) - keep its original "0.00"-style number format but store the
# value as literal text "30.0000", matching the source report export.
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "30.0000"
$ws.Range("P7").NumberFormat = "0.00"

# Number of transactions (Q7) - literal text.
$ws.Range("Q7").NumberFormat = "@"
$ws.Range("Q7").Value = "1:0"

# Totals row - numeric sum of the quantities/prices above.
$ws.Range("P8").Value = 30

# Footer timestamp - bump to the new export time.
$ws.Range("A9").Value = "Wednesday, 8 October, 2025 9:23 AM"
